# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders the "Periodo Mora" rows for both trabajadores (Orlando /
# Alexander) so periods run 2101..2107 in ascending order, interleaving
# each trabajador per period; refreshes "Salario Basico" (col G) to the
# updated value (877803) for every row, and swaps which period (2101 vs
# 2107) carries the smaller "Valor Mora" (col F) amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docType   = "CC"
$orlandoId = "1232597678"
$orlandoNm = "ORLANDO ENRIQUE MENDEZ MORENO"
$alexId    = "1237439959"
$alexNm    = "ALEXANDER HENRIQUE MENDEZ MORENO"

$salario = 877803

# row -> (docType, docId, nombre, periodo, valorMora)
$rows = @{
    16 = @($docType, $orlandoId, $orlandoNm, "2101", 35112)
    17 = @($docType, $alexId,    $alexNm,    "2101", 35112)
    18 = @($docType, $orlandoId, $orlandoNm, "2102", 35112)
    19 = @($docType, $alexId,    $alexNm,    "2102", 35112)
    20 = @($docType, $orlandoId, $orlandoNm, "2103", 35112)
    21 = @($docType, $alexId,    $alexNm,    "2103", 35112)
    22 = @($docType, $orlandoId, $orlandoNm, "2104", 35112)
    23 = @($docType, $alexId,    $alexNm,    "2104", 35112)
    24 = @($docType, $orlandoId, $orlandoNm, "2105", 35112)
    25 = @($docType, $alexId,    $alexNm,    "2105", 35112)
    26 = @($docType, $orlandoId, $orlandoNm, "2106", 35112)
    27 = @($docType, $alexId,    $alexNm,    "2106", 35112)
    28 = @($docType, $orlandoId, $orlandoNm, "2107", 29260)
    29 = @($docType, $alexId,    $alexNm,    "2107", 29260)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $vals[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $salario   # G - Salario Basico
}
